$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore prior revision's value for the "From" column of rule R20 (row 10)
$ws.Range("C10").Value = 1
